# edit.ps1 - apply the "alternate branch" revision to mlk.docx
#
# Changes made:
#  1. Paragraph 1: "This is a Microsoft word document." gets two trailing
#     spaces appended, then three new runs colored C00000 (dark red) are
#     appended spelling out "(This is a change – Version for branch alternate)"
#     (split across three runs, matching the source edit).
#  2. A new, empty paragraph shaded F9F9F9 is appended at the very end of
#     the document body (after the closing "Free at last" paragraph).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the first paragraph's text and append the red annotation.
# ---------------------------------------------------------------------

# Add the two trailing spaces to the existing sentence (keeps it as a
# single run, matching the target markup).
$null = $d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$firstPara = $d.Paragraphs(1)
$insertPoint = $firstPara.Range.End - 1   # just before the paragraph mark

$redColor = 192        # wdColor for RGB(192,0,0) == OOXML w:val="C00000"
$enDash = [char]0x2013  # –

# Run 2: "(This is a change – Ve"
$seg1 = "(This is a change " + $enDash + " Ve"
$r1 = $d.Range($insertPoint, $insertPoint)
$r1.InsertAfter($seg1)
$seg1End = $insertPoint + $seg1.Length
$d.Range($insertPoint, $seg1End).Font.Color = $redColor

# Run 3: "rsion for branch alternate"
$seg2 = "rsion for branch alternate"
$r2 = $d.Range($seg1End, $seg1End)
$r2.InsertAfter($seg2)
$seg2End = $seg1End + $seg2.Length
$d.Range($seg1End, $seg2End).Font.Color = $redColor

# Run 4: ")"
$seg3 = ")"
$r3 = $d.Range($seg2End, $seg2End)
$r3.InsertAfter($seg3)
$seg3End = $seg2End + $seg3.Length
$d.Range($seg2End, $seg3End).Font.Color = $redColor

# ---------------------------------------------------------------------
# 2. Append a new, shaded (F9F9F9), empty paragraph at the end of the
#    document body.
# ---------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Style = "Normal"

$newPara.Range.Shading.Texture = 0                 # wdTextureNone -> shd val="clear"
$newPara.Range.Shading.ForegroundPatternColor = -16777216   # wdColorAutomatic -> color="auto"
$newPara.Range.Shading.BackgroundPatternColor = 16382457    # RGB(0xF9,0xF9,0xF9) -> fill="F9F9F9"

Write-Output "done"
